$wb = $excel.ActiveWorkbook

# Update the "展览" sheet
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 142
$wsExhibit.Range("F4").Value = 96

# Update the "全部类型" sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 142
$wsAll.Range("F4").Value = 96
